$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new job posting row (row 25) with Job_Id=JD_024
$ws.Range("A25").Value = "JD_024"
$ws.Range("B25").Value = "string"
$ws.Range("C25").Value = "string"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
